# Update the Export PDF
# Append two new scan-log rows (37 and 38) to the "Scan Logs" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37
$ws.Range("A37").Value = "2025-11-30T12:31:12.001Z"
$ws.Range("B37").Value = "rus@h.com"
$ws.Range("C37").Value = "https://ksu.edu.sa"
$ws.Range("D37").Value = "Yes"
$ws.Range("E37").Value = "٣٠‏/١١‏/٢٠٢٥"
$ws.Range("F37").Value = "٣:٣١:١٢ م"

# Row 38
$ws.Range("A38").Value = "2025-11-30T12:38:36.960Z"
$ws.Range("B38").Value = "test@test.com"
$ws.Range("C38").Value = "https://ksu.edu.sa"
$ws.Range("D38").Value = "No"
$ws.Range("E38").Value = "٣٠‏/١١‏/٢٠٢٥"
$ws.Range("F38").Value = "٣:٣٨:٣٧ م"
